$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting used by the row above (A8) before filling A9's value,
# since the blank template row used a slightly different border style.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new journal entry on row 9 (Version 0.1 finished - 14/03/2020)
$ws.Range("A9").Value = "Gabriel Pereira"

$entryDate = Get-Date -Year 2020 -Month 3 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B9").Value = $entryDate.Date

$ws.Range("C9").Value = 0.77083333333333337

$ws.Range("D9").Value = "Bataille Navale"
$ws.Range("E9").Value = "Version 0.1 fini"
$ws.Range("F9").Value = "Création de la version 0.1 terminé"

# Update the active selection to match the saved workbook state
$ws.Range("G7").Select()
